$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the two new header cells so that they
# reuse the existing bold/centered/bordered style (s="1") instead of Excel
# creating a brand new style record.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for column I (I0)
$ws.Range("I2").Value = 6
$ws.Range("I3").Value = 4
$ws.Range("I4").Value = 3

# New data values for column J (IF)
$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 4
$ws.Range("J4").Value = 3
